{"js": "// The document has a single table (5 columns x 20 rows) whose cells each\n// contain a short arithmetic expression, e.g. \"53-46=\". This edit replaces\n// every one of the 100 cell values with a new expression, preserving each\n// cell's existing run formatting (font/size) by only touching the text.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst newValues = [\n  [\"96-48=\", \"86-83=\", \"99-65=\", \"64-11=\", \"56+10=\"],\n  [\"84-41=\", \"93-40=\", \"96-2=\", \"99-66=\", \"23+55=\"],\n  [\"38+38=\", \"49-12=\", \"93-60=\", \"60-32=\", \"57-33=\"],\n  [\"79+7=\", \"67-39=\", \"91-42=\", \"89-29=\", \"35-23=\"],\n  [\"7+67=\", \"28-9=\", \"96-2=\", \"74-23=\", \"67-40=\"],\n  [\"27-8=\", \"38-5=\", \"93-60=\", \"65+9=\", \"30-11=\"],\n  [\"30+47=\", \"57+19=\", \"81-50=\", \"86-49=\", \"61-43=\"],\n  [\"67+14=\", \"44+22=\", \"50+8=\", \"27+32=\", \"10+2=\"],\n  [\"56+21=\", \"75-72=\", \"81-47=\", \"96-59=\", \"64-57=\"],\n  [\"93-45=\", \"64-3=\", \"17+52=\", \"79-20=\", \"14+43=\"],\n  [\"0+21=\", \"61-10=\", \"53-6=\", \"61+20=\", \"45-13=\"],\n  [\"89-47=\", \"1+92=\", \"2+12=\", \"63+30=\", \"21+76=\"],\n  [\"39+9=\", \"6+25=\", \"6+27=\", \"44+35=\", \"79-56=\"],\n  [\"98-60=\", \"86-61=\", \"4+38=\", \"94-18=\", \"37-33=\"],\n  [\"66-21=\", \"5+45=\", \"60-6=\", \"36+29=\", \"1+10=\"],\n  [\"38+61=\", \"48+7=\", \"39-4=\", \"27+28=\", \"92-4=\"],\n  [\"68-61=\", \"47-45=\", \"90-34=\", \"36+9=\", \"22+26=\"],\n  [\"68-7=\", \"0+36=\", \"29+2=\", \"69-66=\", \"2+9=\"],\n  [\"72-33=\", \"91-19=\", \"64-52=\", \"28-26=\", \"90-27=\"],\n  [\"82+14=\", \"77+11=\", \"38+21=\", \"99-43=\", \"13+75=\"],\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the 100 arithmetic-problem cells (5 columns x 20 rows) in the\n# first table, in row-major document order, to their new values.\n$d = $word.ActiveDocument\n\n$newValues = @(\n    \"96-48=\",\n    \"86-83=\",\n    \"99-65=\",\n    \"64-11=\",\n    \"56+10=\",\n    \"84-41=\",\n    \"93-40=\",\n    \"96-2=\",\n    \"99-66=\",\n    \"23+55=\",\n    \"38+38=\",\n    \"49-12=\",\n    \"93-60=\",\n    \"60-32=\",\n    \"57-33=\",\n    \"79+7=\",\n    \"67-39=\",\n    \"91-42=\",\n    \"89-29=\",\n    \"35-23=\",\n    \"7+67=\",\n    \"28-9=\",\n    \"96-2=\",\n    \"74-23=\",\n    \"67-40=\",\n    \"27-8=\",\n    \"38-5=\",\n    \"93-60=\",\n    \"65+9=\",\n    \"30-11=\",\n    \"30+47=\",\n    \"57+19=\",\n    \"81-50=\",\n    \"86-49=\",\n    \"61-43=\",\n    \"67+14=\",\n    \"44+22=\",\n    \"50+8=\",\n    \"27+32=\",\n    \"10+2=\",\n    \"56+21=\",\n    \"75-72=\",\n    \"81-47=\",\n    \"96-59=\",\n    \"64-57=\",\n    \"93-45=\",\n    \"64-3=\",\n    \"17+52=\",\n    \"79-20=\",\n    \"14+43=\",\n    \"0+21=\",\n    \"61-10=\",\n    \"53-6=\",\n    \"61+20=\",\n    \"45-13=\",\n    \"89-47=\",\n    \"1+92=\",\n    \"2+12=\",\n    \"63+30=\",\n    \"21+76=\",\n    \"39+9=\",\n    \"6+25=\",\n    \"6+27=\",\n    \"44+35=\",\n    \"79-56=\",\n    \"98-60=\",\n    \"86-61=\",\n    \"4+38=\",\n    \"94-18=\",\n    \"37-33=\",\n    \"66-21=\",\n    \"5+45=\",\n    \"60-6=\",\n    \"36+29=\",\n    \"1+10=\",\n    \"38+61=\",\n    \"48+7=\",\n    \"39-4=\",\n    \"27+28=\",\n    \"92-4=\",\n    \"68-61=\",\n    \"47-45=\",\n    \"90-34=\",\n    \"36+9=\",\n    \"22+26=\",\n    \"68-7=\",\n    \"0+36=\",\n    \"29+2=\",\n    \"69-66=\",\n    \"2+9=\",\n    \"72-33=\",\n    \"91-19=\",\n    \"64-52=\",\n    \"28-26=\",\n    \"90-27=\",\n    \"82+14=\",\n    \"77+11=\",\n    \"38+21=\",\n    \"99-43=\",\n    \"13+75=\"\n)\n\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n\nWrite-Output (\"Updated \" + $i + \" cells\")\n"}
